# feat: add 2022-Q3 data
#
# - Insert a new sheet "2022-Q3" right after "总计" (so the sheet order
#   becomes 总计, 2022-Q3, 2022-Q2, 2021-Q4).
# - Populate "2022-Q3" with the quarter's fund holdings table (8 rows,
#   7 funds) by duplicating the "2022-Q2" sheet (to inherit its header /
#   column-A formatting) and overwriting the values.
# - Update the "总计" summary sheet: the 2022-Q3 totals become the new
#   row 2, and the pre-existing 2022-Q2 / 2021-Q4 rows shift down by one.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2src = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1. Duplicate "2022-Q2" and drop the copy directly after "总计" - this
#    both clones the header/border styling and lands the new tab in its
#    final position in one step (avoids any later Move() call).
# ---------------------------------------------------------------------
$wsQ2src.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Extend the copied sheet from 3 rows (1 header + 2 data) to 8 rows
#    (1 header + 7 data), replicating row 3's formatting down.
# ---------------------------------------------------------------------
$fmtSrc = $wsQ3.Range("A3:H3")
for ($i = 4; $i -le 8; $i++) {
    $dst = $wsQ3.Range("A" + $i + ":H" + $i)
    $fmtSrc.Copy($dst)
}

# ---------------------------------------------------------------------
# 3. Write the 2022-Q3 fund table values (rows 2..8).
#    Col A / H are numbers; B/C/D/E/F/G are text (leading apostrophe
#    keeps numeric-looking strings - fund codes, percentages - as text
#    instead of silently coercing to numbers).
# ---------------------------------------------------------------------
$q3data = @(
    @(0,"001869","招商制造业转型灵活配置混合A","21.64","87.08","3.42","0.7401",9),
    @(1,"213003","宝盈策略增长混合","10.55","90.74","4.95","0.5222",9),
    @(2,"213002","宝盈泛沿海增长混合","5.18","91.39","5.71","0.2958",7),
    @(3,"004569","招商制造业转型灵活配置混合C","6.05","87.08","3.42","0.2069",9),
    @(4,"008132","鹏华价值驱动混合","4.47","92.01","4.09","0.1828",4),
    @(5,"000796","宝盈睿丰创新灵活配置混合 - C","0.64","89.70","5.32","0.0340",8),
    @(6,"000794","宝盈睿丰创新灵活配置混合 - A/B","0.41","89.70","5.32","0.0218",8)
)

$r = 2
foreach ($rec in $q3data) {
    $wsQ3.Cells.Item($r, 1).Value = $rec[0]
    $wsQ3.Cells.Item($r, 2).Value = "'" + $rec[1]
    $wsQ3.Cells.Item($r, 3).Value = $rec[2]
    $wsQ3.Cells.Item($r, 4).Value = "'" + $rec[3]
    $wsQ3.Cells.Item($r, 5).Value = "'" + $rec[4]
    $wsQ3.Cells.Item($r, 6).Value = "'" + $rec[5]
    $wsQ3.Cells.Item($r, 7).Value = "'" + $rec[6]
    $wsQ3.Cells.Item($r, 8).Value = $rec[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. Update the "总计" sheet: shift the existing data rows down by one
#    (2022-Q2 -> row3, 2021-Q4 -> row4) and insert 2022-Q3 as row 2.
#    Row 4 is brand new, so first clone row 3's cell formatting
#    (bordered/bold column A) down into it.
# ---------------------------------------------------------------------
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4:D4"))

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q4"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.09

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.63

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 2
